$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12 (pushes "Application Controller" and
# everything below it down by one row), then populate the new row's
# Element column with "Export Database" (the new tag row for the
# Admin DB export flow). This matches the rest of the sheet, i.e. empty
# produced/sent (B/C/D) cells for this element.
$ws.Rows(12).Insert()
$ws.Range("A12").Value = "Export Database"

Write-Host "Done inserting Export Database row"
